# "Generate Report for Archive"
#
# 1. Update the Status text used throughout the report from
#    "Ready for handoff" to "In Translation" (Overview!E2:F3, and
#    column C rows 2-3 on the zh-cn / de-de sheets).
# 2. Shrink the now-narrower Status column(s) to match the new text:
#    Overview columns E & F, and column C on zh-cn / de-de.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# The Status columns auto-shrank (from a stored width of ~17.22 down to
# ~13.41 "character" units) once the text went from "Ready for handoff"
# (17 chars) to "In Translation" (14 chars). This runtime quantizes
# ColumnWidth to whole pixels (width*6 rounded, plus 5px padding), so
# 12.5 is the input that lands closest to that recorded ~13.41 width.
$newColWidth = 12.5

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $nRows = $used.Rows.Count
    $nCols = $used.Columns.Count
    for ($r = 1; $r -le $nRows; $r++) {
        for ($c = 1; $c -le $nCols; $c++) {
            $cell = $used.Cells.Item($r, $c)
            # Compare as string explicitly: PowerShell's -eq coerces the
            # right-hand side to the left-hand side's type, and some cells
            # hold Booleans, so use the literal string as the left operand.
            if ($oldStatus -eq [string]$cell.Value()) {
                $cell.Value = $newStatus
            }
        }
    }
}

# Overview sheet: zh-cn (E) / de-de (F) status columns
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $newColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColWidth

# zh-cn sheet: Status column (C)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = $newColWidth

# de-de sheet: Status column (C)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = $newColWidth
